$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 3999.8572
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

# Row 23
$ws.Range("H23").Value = 3999.8572
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

# Row 38
$ws.Range("H38").Value = 1925.9584
$ws.Range("I38").Value = 1320.5
$ws.Range("J38").Value = 3136.875
$ws.Range("K38").Value = 3961.5
$ws.Range("L38").Value = 9410.625
$ws.Range("M38").Value = -3589.5
$ws.Range("N38").Value = -10154.625

# Row 51
$ws.Range("H51").Value = 3009.8857
$ws.Range("I51").Value = 3017.8572
$ws.Range("J51").Value = 2997.9285
$ws.Range("K51").Value = 3017.8572
$ws.Range("L51").Value = 2997.9285
$ws.Range("M51").Value = -2533.8572
$ws.Range("N51").Value = -3965.9285

# Row 104
$ws.Range("H104").Value = 460.5
$ws.Range("I104").Value = 460.5
$ws.Range("K104").Value = 1381.5
$ws.Range("M104").Value = 365.5

# Row 138
$ws.Range("H138").Value = 2014.7028
$ws.Range("I138").Value = 870.5714
$ws.Range("J138").Value = 2468.0378
$ws.Range("K138").Value = 2611.7142
$ws.Range("L138").Value = 7404.1134
$ws.Range("M138").Value = 2528.2858
$ws.Range("N138").Value = -17684.1134

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 113394.54
$ws.Range("I32").Value = 121808.37
$ws.Range("K32").Value = 121808.37
$ws.Range("M32").Value = -121521.37

# Row 41
$ws.Range("H41").Value = 2538.5454
$ws.Range("I41").Value = 2042.4
$ws.Range("K41").Value = 2042.4
$ws.Range("M41").Value = -1628.4

# Row 45
$ws.Range("H45").Value = 3732.1538
$ws.Range("I45").Value = 3347.0908
$ws.Range("K45").Value = 3347.0908
$ws.Range("M45").Value = -2970.0908

# Row 61
$ws.Range("H61").Value = 5969.9546
$ws.Range("I61").Value = 1900.1333
$ws.Range("J61").Value = 14691
$ws.Range("K61").Value = 1900.1333
$ws.Range("L61").Value = 14691
$ws.Range("M61").Value = -1688.1333
$ws.Range("N61").Value = -15115

# Row 132
$ws.Range("H132").Value = 2275069
$ws.Range("I132").Value = 2275069
$ws.Range("K132").Value = 6825207
$ws.Range("M132").Value = -6822677

# Row 136
$ws.Range("H136").Value = 5969.9546
$ws.Range("I136").Value = 1900.1333
$ws.Range("J136").Value = 14691
$ws.Range("K136").Value = 5700.3999
$ws.Range("L136").Value = 44073
$ws.Range("M136").Value = -3150.3999
$ws.Range("N136").Value = -49173

$ws = $wb.Worksheets.Item("BSM")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""

# Row 86
$ws.Range("H86").Value = 3192.8235
$ws.Range("I86").Value = 2889.7273
$ws.Range("J86").Value = 3748.5
$ws.Range("K86").Value = 2889.7273
$ws.Range("L86").Value = 3748.5
$ws.Range("M86").Value = -1766.7273
$ws.Range("N86").Value = -5994.5

# Row 89
$ws.Range("H89").Value = 3192.8235
$ws.Range("I89").Value = 2889.7273
$ws.Range("J89").Value = 3748.5
$ws.Range("K89").Value = 14448.6365
$ws.Range("L89").Value = 18742.5
$ws.Range("M89").Value = -8832.636500000001
$ws.Range("N89").Value = -29974.5

# Row 97
$ws.Range("H97").Value = 38989.6
$ws.Range("I97").Value = 9999.5
$ws.Range("K97").Value = 9999.5
$ws.Range("M97").Value = -9008.5

# Row 105
$ws.Range("H105").Value = 2425.3635
$ws.Range("I105").Value = 1962.425
$ws.Range("K105").Value = 1962.425
$ws.Range("M105").Value = -215.425

# Row 107
$ws.Range("H107").Value = 652.2593000000001
$ws.Range("I107").Value = 624.52
$ws.Range("K107").Value = 624.52
$ws.Range("M107").Value = 1295.48

# Row 134
$ws.Range("H134").Value = 22802.2
$ws.Range("I134").Value = 2006
$ws.Range("K134").Value = 6018
$ws.Range("M134").Value = -3483

$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Range("H21").Value = 2505.5
$ws.Range("I21").Value = 2505.5
$ws.Range("K21").Value = 2505.5
$ws.Range("M21").Value = -2270.5

# Row 44
$ws.Range("H44").Value = 1364.75
$ws.Range("I44").Value = 1364.75
$ws.Range("K44").Value = 1364.75
$ws.Range("M44").Value = -922.75

# Row 55
$ws.Range("H55").Value = 8111
$ws.Range("I55").Value = 8111
$ws.Range("K55").Value = 8111
$ws.Range("M55").Value = -7796

# Row 94
$ws.Range("H94").Value = 5683.7393
$ws.Range("I94").Value = 10354.182
$ws.Range("K94").Value = 10354.182
$ws.Range("M94").Value = -9903.182000000001

# Row 132
$ws.Range("H132").Value = 2463.7646
$ws.Range("I132").Value = 2182.3076
$ws.Range("K132").Value = 6546.9228
$ws.Range("M132").Value = -4016.9228

# Row 134
$ws.Range("H134").Value = 2187.1875
$ws.Range("I134").Value = 2044.2858
$ws.Range("J134").Value = 3187.5
$ws.Range("K134").Value = 6132.857400000001
$ws.Range("L134").Value = 9562.5
$ws.Range("M134").Value = -3597.857400000001
$ws.Range("N134").Value = -14632.5

$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Range("H24").Value = 271.25
$ws.Range("I24").Value = 217.5
$ws.Range("J24").Value = 325
$ws.Range("K24").Value = 652.5
$ws.Range("L24").Value = 975
$ws.Range("M24").Value = -422.5
$ws.Range("N24").Value = -1435

# Row 62
$ws.Range("H62").Value = 1949.5
$ws.Range("I62").Value = 899
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2697
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -2011
$ws.Range("N62").Value = -10372

# Row 65
$ws.Range("H65").Value = 1949.5
$ws.Range("I65").Value = 899
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 8091
$ws.Range("L65").Value = 27000
$ws.Range("M65").Value = -4659
$ws.Range("N65").Value = -33864

# Row 107
$ws.Range("H107").Value = 1353.4736
$ws.Range("J107").Value = 1640.5385
$ws.Range("L107").Value = 4921.6155
$ws.Range("N107").Value = -8761.6155

# Row 123
$ws.Range("H123").Value = 11000.777
$ws.Range("J123").Value = 14499.25
$ws.Range("L123").Value = 43497.75
$ws.Range("N123").Value = -48397.75

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 42250
$ws.Range("J15").Value = 42250
$ws.Range("L15").Value = 42250
$ws.Range("N15").Value = -42826

# Row 18
$ws.Range("H18").Value = 4739.75
$ws.Range("I18").Value = 4739.75
$ws.Range("K18").Value = 4739.75
$ws.Range("M18").Value = -4446.75

# Row 42
$ws.Range("H42").Value = 100000
$ws.Range("J42").Value = 100000
$ws.Range("L42").Value = 100000
$ws.Range("N42").Value = -100970

# Row 49
$ws.Range("H49").Value = 20025.5
$ws.Range("I49").Value = 20025.5
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 20025.5
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -19841.5
$ws.Range("N49").Value = ""

# Row 81
$ws.Range("H81").Value = 42250
$ws.Range("J81").Value = 42250
$ws.Range("L81").Value = 42250
$ws.Range("N81").Value = -44246

# Row 84
$ws.Range("H84").Value = 42250
$ws.Range("J84").Value = 42250
$ws.Range("L84").Value = 126750
$ws.Range("N84").Value = -136734

# Row 115
$ws.Range("H115").Value = 100000
$ws.Range("J115").Value = 100000
$ws.Range("L115").Value = 100000
$ws.Range("N115").Value = -102350

# Row 122
$ws.Range("H122").Value = 53461.15
$ws.Range("I122").Value = 94239.09
$ws.Range("J122").Value = 3621.4443
$ws.Range("K122").Value = 282717.27
$ws.Range("L122").Value = 10864.3329
$ws.Range("M122").Value = -280267.27
$ws.Range("N122").Value = -15764.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""

# Row 16
$ws.Range("H16").Value = 369
$ws.Range("I16").Value = 369
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 369
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -199
$ws.Range("N16").Value = ""

# Row 22
$ws.Range("H22").Value = 2093.8857
$ws.Range("I22").Value = 599
$ws.Range("J22").Value = 2234.0312
$ws.Range("K22").Value = 599
$ws.Range("L22").Value = 2234.0312
$ws.Range("M22").Value = -304
$ws.Range("N22").Value = -2824.0312

# Row 27
$ws.Range("H27").Value = 2093.8857
$ws.Range("I27").Value = 599
$ws.Range("J27").Value = 2234.0312
$ws.Range("K27").Value = 599
$ws.Range("L27").Value = 2234.0312
$ws.Range("M27").Value = -492
$ws.Range("N27").Value = -2448.0312

# Row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""

# Row 55
$ws.Range("H55").Value = 2018.1072
$ws.Range("I55").Value = 2461.625
$ws.Range("J55").Value = 1840.7
$ws.Range("K55").Value = 2461.625
$ws.Range("L55").Value = 1840.7
$ws.Range("M55").Value = -2288.625
$ws.Range("N55").Value = -2186.7

# Row 61
$ws.Range("H61").Value = 11911.053
$ws.Range("I61").Value = 11380.143
$ws.Range("K61").Value = 11380.143
$ws.Range("M61").Value = -11178.143

# Row 113
$ws.Range("H113").Value = 11911.053
$ws.Range("I113").Value = 11380.143
$ws.Range("K113").Value = 11380.143
$ws.Range("M113").Value = -9210.143

# Row 136
$ws.Range("H136").Value = 10496.9375
$ws.Range("I136").Value = 4596.1
$ws.Range("K136").Value = 13788.3
$ws.Range("M136").Value = -11238.3

$ws = $wb.Worksheets.Item("WVR")
# Row 115
$ws.Range("H115").Value = 49998
$ws.Range("J115").Value = 49998
$ws.Range("L115").Value = 49998
$ws.Range("N115").Value = -53132

# Row 136
$ws.Range("H136").Value = 4199.909
$ws.Range("I136").Value = 3783.3333
$ws.Range("K136").Value = 11349.9999
$ws.Range("M136").Value = -8799.999899999999
